$wb = $excel.ActiveWorkbook

# --- Sheet "Valeurs reelles": update S+1/S+2/S+3 headers to *_class and rewrite classification grid ---
$ws1 = $wb.Worksheets.Item("Valeurs réelles")

$ws1.Range("C1").Value = "PRIX EXP POMME GOLDEN FRANCE 115/150G CAT.I SACHET_S+1_class"
$ws1.Range("D1").Value = "PRIX EXP POMME GOLDEN FRANCE 115/150G CAT.I SACHET_S+2_class"
$ws1.Range("E1").Value = "PRIX EXP POMME GOLDEN FRANCE 115/150G CAT.I SACHET_S+3_class"

$sheet1Grid = @(
    @(4, 2, 2),
    @(2, 2, 2),
    @(2, 2, 2),
    @(2, 2, 2),
    @(2, 2, 2),
    @(2, 2, 2),
    @(2, 2, 2),
    @(2, 2, 4),
    @(2, 4, 2),
    @(4, 2, 2),
    @(2, 2, 1),
    @(2, 1, 2),
    @(1, 2, 1),
    @(2, 1, 2),
    @(1, 2, 1),
    @(2, 1, 2),
    @(1, 2, 0),
    @(2, 0, 3),
    @(0, 3, 1),
    @(3, 1, 3),
    @(1, 3, 1),
    @(3, 1, 3),
    @(1, 3, 0),
    @(3, 0, 1),
    @(0, 1, 2),
    @(1, 2, 2),
    @(2, 2, 2)
)

for ($i = 0; $i -lt $sheet1Grid.Length; $i++) {
    $row = 2 + $i
    $ws1.Cells.Item($row, 3).Value = $sheet1Grid[$i][0]
    $ws1.Cells.Item($row, 4).Value = $sheet1Grid[$i][1]
    $ws1.Cells.Item($row, 5).Value = $sheet1Grid[$i][2]
}

# --- Sheet "Predictions": rewrite PRED_S1/S2/S3 grid with classifier output ---
$ws2 = $wb.Worksheets.Item("Prédictions")

$sheet2Grid = @(
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 2, 2),
    @(0, 0, 0),
    @(0, 2, 2),
    @(0, 2, 0),
    @(2, 2, -2),
    @(0, 0, 0),
    @(2, 2, 0),
    @(2, 0, 0),
    @(-2, 0, 0),
    @(0, 0, -2),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(-2, -2, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(-2, 0, 0),
    @(-2, 0, 0),
    @(0, 0, 0),
    @(-1, 0, -2),
    @(0, -2, 2)
)

for ($i = 0; $i -lt $sheet2Grid.Length; $i++) {
    $row = 2 + $i
    $ws2.Cells.Item($row, 2).Value = $sheet2Grid[$i][0]
    $ws2.Cells.Item($row, 3).Value = $sheet2Grid[$i][1]
    $ws2.Cells.Item($row, 4).Value = $sheet2Grid[$i][2]
}
